$wb = $excel.ActiveWorkbook

# Set the boolean control-lever value on the BIEfIE sheet (B2) from 1 to 0
$ws = $wb.Worksheets.Item("BIEfIE")
$ws.Range("B2").Value = 0

# Make "About" the active sheet/tab (matches tabSelected moving from BIEfIE to About)
$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Activate()
